$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Nicolás"
$ws.Range("C2").Value = "nicolas72lol@gmail.com"
$ws.Range("D2").Value = "8c6976e5b5410415bde908bd4dee15dfb167a9c873fc4bb8a81f6f2ab448a918"
$ws.Range("E2").Value = "admin"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Usuario"
$ws.Range("C3").Value = "pepe@gmail.com"
$ws.Range("D3").Value = "03ac674216f3e15c761ee1a5e255f067953623c8b388b4459e13f978d7c846f4"
$ws.Range("E3").Value = "usuario"
